$d = $word.ActiveDocument

$replacements = @(
    @("398×8=", "169×9="),
    @("493×5=", "314×9="),
    @("300×3=", "344×7="),
    @("488×3=", "386×3="),
    @("979×9=", "263×8="),
    @("743×7=", "960×6="),
    @("847×9=", "805×2="),
    @("910×8=", "757×7="),
    @("898×3=", "791×4="),
    @("670×6=", "356×8="),
    @("553×3=", "944×7="),
    @("629×8=", "517×9="),
    @("320×3=", "149×9="),
    @("935×5=", "556×8="),
    @("643×2=", "983×6="),
    @("175×6=", "713×4="),
    @("590×8=", "994×4="),
    @("566×9=", "757×9="),
    @("712×3=", "295×4="),
    @("314×2=", "359×3="),
    @("950×7=", "246×2="),
    @("683×3=", "410×5="),
    @("700×2=", "698×5="),
    @("951×4=", "295×7="),
    @("281×8=", "190×5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
